# Target the worksheet named "Card1" (the edit applies to sheet Card1, per the commit message)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# Row 3, columns E..L were blank inline-string cells; they become the text "nan"
$ws.Range("E3:L3").Value = "nan"

# Column P on row 3 was also blank; it becomes "nan" as well
$ws.Range("P3").Value = "nan"
